$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2:E2").NumberFormat = "@"
$ws.Range("D2").Value = '26.608.06'
$ws.Range("E2").Value = '  +1.19%  '
$ws.Range("D2:E2").ClearFormats()

# Row 3
$ws.Range("D3:E3").NumberFormat = "@"
$ws.Range("D3").Value = '1.630.91'
$ws.Range("E3").Value = '  +1.39%  '
$ws.Range("D3:E3").ClearFormats()

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("E4").ClearFormats()

# Row 5
$ws.Range("D5:E5").NumberFormat = "@"
$ws.Range("D5").Value = '212.59'
$ws.Range("E5").Value = '  -0.15%  '
$ws.Range("D5:E5").ClearFormats()

# Row 6
$ws.Range("B6:E6").NumberFormat = "@"
$ws.Range("B6").Value = 'USDC'
$ws.Range("C6").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("D6").Value = '1.00'
$ws.Range("E6").Value = '  +0.07%  '
$ws.Range("B6:E6").ClearFormats()

# Row 7
$ws.Range("B7:E7").NumberFormat = "@"
$ws.Range("B7").Value = 'XRP'
$ws.Range("C7").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D7").Value = '0.493'
$ws.Range("E7").Value = '  +1.40%  '
$ws.Range("B7:E7").ClearFormats()

# Row 8
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.75%  '
$ws.Range("E8").ClearFormats()

# Row 9
$ws.Range("D9:E9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0623'
$ws.Range("E9").Value = '  +1.40%  '
$ws.Range("D9:E9").ClearFormats()

# Row 10
$ws.Range("D10:E10").NumberFormat = "@"
$ws.Range("D10").Value = '18.94'
$ws.Range("E10").Value = '  +2.80%  '
$ws.Range("D10:E10").ClearFormats()

# Row 11
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +3.38%  '
$ws.Range("E11").ClearFormats()

# Row 12
$ws.Range("D12:E12").NumberFormat = "@"
$ws.Range("D12").Value = '1.859.29'
$ws.Range("E12").Value = '  +1.45%  '
$ws.Range("D12:E12").ClearFormats()

# Row 13
$ws.Range("D13:E13").NumberFormat = "@"
$ws.Range("D13").Value = '1.621.42'
$ws.Range("E13").Value = '  +0.87%  '
$ws.Range("D13:E13").ClearFormats()

# Row 14
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +1.22%  '
$ws.Range("E14").ClearFormats()

# Row 15
$ws.Range("D15:E15").NumberFormat = "@"
$ws.Range("D15").Value = '0.526'
$ws.Range("E15").Value = '  +1.99%  '
$ws.Range("D15:E15").ClearFormats()

# Row 16
$ws.Range("D16:E16").NumberFormat = "@"
$ws.Range("D16").Value = '26.600.34'
$ws.Range("E16").Value = '  +1.28%  '
$ws.Range("D16:E16").ClearFormats()

# Row 17
$ws.Range("D17:E17").NumberFormat = "@"
$ws.Range("D17").Value = '62.95'
$ws.Range("E17").Value = '  +1.26%  '
$ws.Range("D17:E17").ClearFormats()

# Row 18
$ws.Range("D18:E18").NumberFormat = "@"
$ws.Range("D18").Value = '0.0₃0738'
$ws.Range("E18").Value = '  +1.47%  '
$ws.Range("D18:E18").ClearFormats()

# Row 19
$ws.Range("D19:E19").NumberFormat = "@"
$ws.Range("D19").Value = '209.28'
$ws.Range("E19").Value = '  +3.77%  '
$ws.Range("D19:E19").ClearFormats()

# Row 20
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.04%  '
$ws.Range("E20").ClearFormats()

# Row 21
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +0.66%  '
$ws.Range("E21").ClearFormats()

# Row 22
$ws.Range("D22:E22").NumberFormat = "@"
$ws.Range("D22").Value = '9.41'
$ws.Range("E22").Value = '  +0.83%  '
$ws.Range("D22:E22").ClearFormats()

# Row 23
$ws.Range("D23:E23").NumberFormat = "@"
$ws.Range("D23").Value = '6.17'
$ws.Range("E23").Value = '  +2.63%  '
$ws.Range("D23:E23").ClearFormats()

# Row 24
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +2.92%  '
$ws.Range("E24").ClearFormats()

# Row 25
$ws.Range("D25:E25").NumberFormat = "@"
$ws.Range("D25").Value = '146.95'
$ws.Range("E25").Value = '  +2.56%  '
$ws.Range("D25:E25").ClearFormats()

# Row 26
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.01%  '
$ws.Range("E26").ClearFormats()

# Row 27
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -0.58%  '
$ws.Range("E27").ClearFormats()

# Row 28
$ws.Range("D28:E28").NumberFormat = "@"
$ws.Range("D28").Value = '6.84'
$ws.Range("E28").Value = '  +4.11%  '
$ws.Range("D28:E28").ClearFormats()

# Row 29
$ws.Range("D29:E29").NumberFormat = "@"
$ws.Range("D29").Value = '15.34'
$ws.Range("E29").Value = '  +0.72%  '
$ws.Range("D29:E29").ClearFormats()

# Row 30
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +4.86%  '
$ws.Range("E30").ClearFormats()

# Row 31
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -0.47%  '
$ws.Range("E31").ClearFormats()

# Row 33
$ws.Range("D33:E33").NumberFormat = "@"
$ws.Range("D33").Value = '2.94'
$ws.Range("E33").Value = '  -0.07%  '
$ws.Range("D33:E33").ClearFormats()

# Row 34
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +0.45%  '
$ws.Range("E34").ClearFormats()

# Row 35
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -0.02%  '
$ws.Range("E35").ClearFormats()

# Row 36
$ws.Range("D36:E36").NumberFormat = "@"
$ws.Range("D36").Value = '1.165.02'
$ws.Range("E36").Value = '  +0.33%  '
$ws.Range("D36:E36").ClearFormats()

# Row 37
$ws.Range("D37:E37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0169'
$ws.Range("E37").Value = '  +0.26%  '
$ws.Range("D37:E37").ClearFormats()

# Row 38
$ws.Range("D38:E38").NumberFormat = "@"
$ws.Range("D38").Value = '0.805'
$ws.Range("E38").Value = '  +2.20%  '
$ws.Range("D38:E38").ClearFormats()

# Row 39
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +0.06%  '
$ws.Range("E39").ClearFormats()

# Row 40
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +1.47%  '
$ws.Range("E40").ClearFormats()

# Row 41
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -0.24%  '
$ws.Range("E41").ClearFormats()

# Row 42
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +0.99%  '
$ws.Range("E42").ClearFormats()

# Row 43
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +0.05%  '
$ws.Range("E43").ClearFormats()

# Row 44
$ws.Range("D44:E44").NumberFormat = "@"
$ws.Range("D44").Value = '1.770.39'
$ws.Range("E44").Value = '  +1.53%  '
$ws.Range("D44:E44").ClearFormats()

# Row 45
$ws.Range("D45:E45").NumberFormat = "@"
$ws.Range("D45").Value = '91.98'
$ws.Range("E45").Value = '  -0.01%  '
$ws.Range("D45:E45").ClearFormats()

# Row 46
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +0.79%  '
$ws.Range("E46").ClearFormats()

# Row 47
$ws.Range("D47:E47").NumberFormat = "@"
$ws.Range("D47").Value = '54.56'
$ws.Range("E47").Value = '  +0.96%  '
$ws.Range("D47:E47").ClearFormats()

# Row 48
$ws.Range("B48:E48").NumberFormat = "@"
$ws.Range("B48").Value = 'BabyDogeCoin'
$ws.Range("C48").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D48").Value = '0.0₆0103'
$ws.Range("E48").Value = '  -2.90%  '
$ws.Range("B48:E48").ClearFormats()

# Row 49
$ws.Range("B49:E49").NumberFormat = "@"
$ws.Range("B49").Value = 'Cronos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D49").Value = '0.0510'
$ws.Range("E49").Value = '  +0.72%  '
$ws.Range("B49:E49").ClearFormats()

# Row 50
$ws.Range("D50:E50").NumberFormat = "@"
$ws.Range("D50").Value = '7.55'
$ws.Range("E50").Value = '  +4.10%  '
$ws.Range("D50:E50").ClearFormats()

# Row 51
$ws.Range("B51:E51").NumberFormat = "@"
$ws.Range("B51").Value = 'Mantle'
$ws.Range("C51").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D51").Value = '0.409'
$ws.Range("E51").Value = '  +0.49%  '
$ws.Range("B51:E51").ClearFormats()
